# Update the date heading.
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Find.Execute("2026-02-27 Friday", $true, $false, $false, $false, $false, `
    $true, 1, $false, "2026-02-28 Saturday", 2)

# Update the division-problem table. Only the text inside each cell changes;
# the table's row/column layout (5 columns, data on rows 1/5/9/13/17) is unchanged.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "92÷2="
$t.Cell(1, 2).Range.Text = "54÷5="
$t.Cell(1, 3).Range.Text = "68÷9="
$t.Cell(1, 4).Range.Text = "25÷9="
$t.Cell(1, 5).Range.Text = "84÷3="

$t.Cell(5, 1).Range.Text = "35÷4="
$t.Cell(5, 2).Range.Text = "52÷9="
$t.Cell(5, 3).Range.Text = "69÷7="
$t.Cell(5, 4).Range.Text = "76÷6="
$t.Cell(5, 5).Range.Text = "59÷2="

$t.Cell(9, 1).Range.Text = "52÷5="
$t.Cell(9, 2).Range.Text = "42÷8="
$t.Cell(9, 3).Range.Text = "18÷5="
$t.Cell(9, 4).Range.Text = "54÷7="
$t.Cell(9, 5).Range.Text = "18÷7="

$t.Cell(13, 1).Range.Text = "68÷6="
$t.Cell(13, 2).Range.Text = "72÷6="
$t.Cell(13, 3).Range.Text = "36÷3="
$t.Cell(13, 4).Range.Text = "68÷9="
$t.Cell(13, 5).Range.Text = "97÷5="

$t.Cell(17, 1).Range.Text = "97÷6="
$t.Cell(17, 2).Range.Text = "18÷6="
$t.Cell(17, 3).Range.Text = "79÷3="
$t.Cell(17, 4).Range.Text = "70÷4="
$t.Cell(17, 5).Range.Text = "86÷9="
